$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated symbol list values (price/volume/hora) per commit diff.
# NumberFormat is forced to Text ("@") before assignment so that
# numeric-looking strings (prices, percentages) are stored as text,
# matching the original inline-string cell contents.

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "302.92"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "1.27%"
$ws.Range("G2").NumberFormat = "@"
$ws.Range("G2").Value = "6"
# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "32.80"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "4.60%"
$ws.Range("G3").NumberFormat = "@"
$ws.Range("G3").Value = "6"
# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "4.955"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "-3.16%"
$ws.Range("G4").NumberFormat = "@"
$ws.Range("G4").Value = "6"
# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.07829"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "-1.47%"
$ws.Range("G5").NumberFormat = "@"
$ws.Range("G5").Value = "6"
# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.965"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "-14.95%"
$ws.Range("G6").NumberFormat = "@"
$ws.Range("G6").Value = "6"
# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "7.879"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "1.03%"
$ws.Range("G7").NumberFormat = "@"
$ws.Range("G7").Value = "6"
# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.807"
$ws.Range("G8").NumberFormat = "@"
$ws.Range("G8").Value = "6"
# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.9280"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "0.61%"
$ws.Range("G9").NumberFormat = "@"
$ws.Range("G9").Value = "6"
# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1763"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "1.38%"
$ws.Range("G10").NumberFormat = "@"
$ws.Range("G10").Value = "6"
# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07787"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "3.43%"
$ws.Range("G11").NumberFormat = "@"
$ws.Range("G11").Value = "6"
# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08661"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-7.84%"
$ws.Range("G12").NumberFormat = "@"
$ws.Range("G12").Value = "6"
# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.03141"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "4.44%"
$ws.Range("G13").NumberFormat = "@"
$ws.Range("G13").Value = "6"
# Row 14
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "0.12%"
$ws.Range("G14").NumberFormat = "@"
$ws.Range("G14").Value = "6"
# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.001510"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "0.21%"
$ws.Range("G15").NumberFormat = "@"
$ws.Range("G15").Value = "6"
# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.005738"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "-6.67%"
$ws.Range("G16").NumberFormat = "@"
$ws.Range("G16").Value = "6"
# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.464"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-0.28%"
$ws.Range("G17").NumberFormat = "@"
$ws.Range("G17").Value = "6"
# Row 18
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "-4.83%"
$ws.Range("G18").NumberFormat = "@"
$ws.Range("G18").Value = "6"
# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.3310"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "1.19%"
$ws.Range("G19").NumberFormat = "@"
$ws.Range("G19").Value = "6"
# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.1318"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "0.77%"
$ws.Range("G20").NumberFormat = "@"
$ws.Range("G20").Value = "6"
# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.295"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "10.14%"
$ws.Range("G21").NumberFormat = "@"
$ws.Range("G21").Value = "6"
# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.1991"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "17.10%"
$ws.Range("G22").NumberFormat = "@"
$ws.Range("G22").Value = "6"
# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04560"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "-1.45%"
$ws.Range("G23").NumberFormat = "@"
$ws.Range("G23").Value = "6"
# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001225"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "-2.22%"
$ws.Range("G24").NumberFormat = "@"
$ws.Range("G24").Value = "6"
# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004442"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "-0.91%"
$ws.Range("G25").NumberFormat = "@"
$ws.Range("G25").Value = "6"
# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0001251"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "4.16%"
$ws.Range("G26").NumberFormat = "@"
$ws.Range("G26").Value = "6"
# Row 27
$ws.Range("G27").NumberFormat = "@"
$ws.Range("G27").Value = "6"
# Row 28
$ws.Range("G28").NumberFormat = "@"
$ws.Range("G28").Value = "6"
# Row 29
$ws.Range("G29").NumberFormat = "@"
$ws.Range("G29").Value = "6"
# Row 30
$ws.Range("G30").NumberFormat = "@"
$ws.Range("G30").Value = "6"
# Row 31
$ws.Range("G31").NumberFormat = "@"
$ws.Range("G31").Value = "6"
# Row 32
$ws.Range("G32").NumberFormat = "@"
$ws.Range("G32").Value = "6"
# Row 33
$ws.Range("G33").NumberFormat = "@"
$ws.Range("G33").Value = "6"
# Row 34
$ws.Range("G34").NumberFormat = "@"
$ws.Range("G34").Value = "6"
# Row 35
$ws.Range("G35").NumberFormat = "@"
$ws.Range("G35").Value = "6"
# Row 36
$ws.Range("G36").NumberFormat = "@"
$ws.Range("G36").Value = "6"
# Row 37
$ws.Range("G37").NumberFormat = "@"
$ws.Range("G37").Value = "6"
# Row 38
$ws.Range("G38").NumberFormat = "@"
$ws.Range("G38").Value = "6"
# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01719"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "-1.80%"
$ws.Range("G39").NumberFormat = "@"
$ws.Range("G39").Value = "6"
# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04775"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "3.35%"
$ws.Range("G40").NumberFormat = "@"
$ws.Range("G40").Value = "6"
# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007492"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "7.67%"
$ws.Range("G41").NumberFormat = "@"
$ws.Range("G41").Value = "6"
# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1357"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "-0.30%"
$ws.Range("G42").NumberFormat = "@"
$ws.Range("G42").Value = "6"
# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.002341"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "6.84%"
$ws.Range("G43").NumberFormat = "@"
$ws.Range("G43").Value = "6"
# Row 44
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "14.20%"
$ws.Range("G44").NumberFormat = "@"
$ws.Range("G44").Value = "6"
# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00006250"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "-0.24%"
$ws.Range("G45").NumberFormat = "@"
$ws.Range("G45").Value = "6"
# Row 46
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "0.09%"
$ws.Range("G46").NumberFormat = "@"
$ws.Range("G46").Value = "6"
# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.003102"
$ws.Range("G47").NumberFormat = "@"
$ws.Range("G47").Value = "6"
# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.8234"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "10.28%"
$ws.Range("G48").NumberFormat = "@"
$ws.Range("G48").Value = "6"
# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.00002101"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "0.09%"
$ws.Range("G49").NumberFormat = "@"
$ws.Range("G49").Value = "6"
# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0002001"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "0.09%"
$ws.Range("G50").NumberFormat = "@"
$ws.Range("G50").Value = "6"
# Row 51
$ws.Range("G51").NumberFormat = "@"
$ws.Range("G51").Value = "6"
